$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$hf = $master.HeadersFooters.DateAndTime
Write-Output "Format before=$($hf.Format)"
$hf.Format = 2
Write-Output "Format after=$($hf.Format)"
Write-Output "Text=$($hf.Text)"
